$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Iteration 2 (row 3) metrics ---
$ws.Range("B3").Value = 165.41
$ws.Range("C3").Value = 17.66
$ws.Range("D3").Value = 6.5
$ws.Range("E3").Value = 12100.01
$ws.Range("F3").Value = -1.56

# --- Update Iteration 3 (row 4) metrics ---
$ws.Range("B4").Value = 293.02
$ws.Range("C4").Value = 30.49
$ws.Range("D4").Value = 25.98
$ws.Range("E4").Value = 8369.94
$ws.Range("F4").Value = 1172.07

# --- Update Iteration 4 (row 5) metrics ---
$ws.Range("B5").Value = 132.44999999999999
$ws.Range("C5").Value = 14.05
$ws.Range("D5").Value = 2.4
$ws.Range("E5").Value = 8083.11
$ws.Range("F5").Value = 112.67

# --- Notes: boundary-penalty note replaced by a re-verify note ---
$ws.Range("G4").Value = "Need to re-verify Num2Move"
$ws.Range("G5").Value = "Need to re-verify Num2Move"

# --- Apply a light background fill across the whole used area ---
$ws.Range("A1:G7").Interior.ThemeColor = 2

# --- Italicize the iteration numbers that were re-run (Iters 2 & 3) ---
$ws.Range("A3:A4").Font.Italic = $true

# --- Print setup tweak ---
$ws.PageSetup.Orientation = 1

# --- Selection state ---
$ws.Range("G15").Select()
